# ---------------------------------------------------------------------------
# Populate the "min" sheet with the full per-category minimum-staffing counts
# (columns C:I) and a running total in column K, then build a new "diff"
# sheet that shows current - min for every category/column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "min" worksheet - fill in the detail columns (C:I) row by row
# ------------------------------------------------------------------
$min = $wb.Worksheets.Item("min")

# row 2 - qrl / priority 1
$min.Range("C2").Value = 0
$min.Range("D2").Value = 0
$min.Range("E2").Value = 5
$min.Range("F2").Value = 0
$min.Range("G2").Value = 1
$min.Range("H2").Value = 5
$min.Range("I2").Value = 1

# row 3 - payroll / priority 1
$min.Range("F3").Value = 7
$min.Range("I3").Value = 2

# row 4 - routine / priority 2
$min.Range("C4").Value = 72
$min.Range("E4").Value = 148
$min.Range("F4").Value = 4
$min.Range("G4").Value = 0
$min.Range("H4").Value = 0
$min.Range("I4").Value = 16

# row 5 - nwts / priority 2
$min.Range("C5").Value = 0
$min.Range("D5").Value = 6
$min.Range("E5").Value = 4
$min.Range("F5").Value = 0
$min.Range("G5").Value = 1
$min.Range("H5").Value = 1
$min.Range("I5").Value = 1

# row 6 - office / priority 3
$min.Range("F6").Value = 2

# row 7 - conv_ctr / priority 3
$min.Range("C7").Value = 6
$min.Range("E7").Value = 6
$min.Range("F7").Value = 1
$min.Range("I7").Value = 1

# row 8 - ds_da / priority 3
$min.Range("C8").Value = 4
$min.Range("E8").Value = 8
$min.Range("F8").Value = 1
$min.Range("I8").Value = 1

# row 9 - mss / priority 3
$min.Range("C9").Value = 2

# row 10 - prop / priority 4 (no detail counts)

# row 11 - spec / priority 4
$min.Range("C11").Value = 5
$min.Range("E11").Value = 10
$min.Range("F11").Value = 1
$min.Range("I11").Value = 1

# row 12 - marine / priority 4
$min.Range("H12").Value = 2

# row 13 - night / priority 4
$min.Range("C13").Value = 1
$min.Range("E13").Value = 2
$min.Range("I13").Value = 1

# Running total column (K) - SUM across C:I for every data row
$min.Range("K2").Formula = "=SUM(C2:I2)"
$min.Range("K3:K13").FormulaR1C1 = "=SUM(RC[-8]:RC[-2])"

# Column J ("emer_drivers") is wide in the other sheets - match that on "min"
$min.Columns.Item(10).AutoFit()

# Restore the selection Excel had when the workbook was last saved
[void]$min.Range("C4").Select()

# ------------------------------------------------------------------
# 2) New "diff" worksheet - current minus min, for every column
# ------------------------------------------------------------------
$current = $wb.Worksheets.Item("current")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$diff = $wb.Worksheets.Add($null, $lastSheet)
$diff.Name = "diff"

# Header row - same labels as "current" / "min", minus the min_req column
$diff.Range("A1").Value = $current.Range("A1").Value2
$diff.Range("B1").Value = $current.Range("B1").Value2
$diff.Range("C1").Value = $current.Range("C1").Value2
$diff.Range("D1").Value = $current.Range("D1").Value2
$diff.Range("E1").Value = $current.Range("E1").Value2
$diff.Range("F1").Value = $current.Range("F1").Value2
$diff.Range("G1").Value = $current.Range("G1").Value2
$diff.Range("H1").Value = $current.Range("H1").Value2
$diff.Range("I1").Value = $current.Range("I1").Value2
$diff.Range("J1").Value = $current.Range("J1").Value2

for ($r = 2; $r -le 13; $r++) {
    $diff.Range("A$r").Value = $current.Range("A$r").Value2
    $diff.Range("B$r").Value = $current.Range("B$r").Value2
}

# Body - current!<cell> - min!<cell> for every data row/column
$diff.Range("C2:J13").FormulaR1C1 = "=current!RC-min!RC"

# Same column-J width treatment as on "min"
$diff.Columns.Item(10).AutoFit()

[void]$diff.Range("I18").Select()

# ------------------------------------------------------------------
# 3) Cosmetic bits on the "current" sheet (selection moved before save)
# ------------------------------------------------------------------
[void]$current.Range("L7").Select()
$current.Activate()
